$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -8.5485798740429075
$ws.Range("D2").Value = -7.8186954769860115
$ws.Range("E2").Value = 76.115569507737561
$ws.Range("F2").Value = 1.0637284526118578
$ws.Range("G2").Value = 3.7012699802074374
$ws.Range("H2").Value = 6.0933217976940091
$ws.Range("I2").Value = 6.0503894393303073
$ws.Range("J2").Value = 2.364772527758932
$ws.Range("K2").Value = 2.0302092899651214
$ws.Range("L2").Value = 1.9419944169125303
$ws.Range("M2").Value = 2.1530192012362019
$ws.Range("N2").Value = 2.090847134109556
$ws.Range("O2").Value = 2.1658217118195751
$ws.Range("P2").Value = 1.8761890347670953
$ws.Range("Q2").Value = 1.5393360456542964
$ws.Range("R2").Value = 2.3333275224678585
$ws.Range("S2").Value = 1.5054315570143899
$ws.Range("T2").Value = 1.97146598328208
$ws.Range("U2").Value = 0.56268936636745548
$ws.Range("V2").Value = 0.6502012509918863
$ws.Range("W2").Value = 0.37607084722408091
$ws.Range("X2").Value = 0.040314128773576173

# Row 3
$ws.Range("C3").Value = -6.3772198878785193
$ws.Range("D3").Value = -7.1241338228082833
$ws.Range("E3").Value = 51.26140278009035
$ws.Range("F3").Value = 0.60674509297341872
$ws.Range("G3").Value = 1.8198975576403746
$ws.Range("H3").Value = 4.2944577880489367
$ws.Range("I3").Value = 10.517490430560629
$ws.Range("J3").Value = 3.9581710288706606
$ws.Range("K3").Value = 2.8901025888258571
$ws.Range("L3").Value = 3.3218757092161066
$ws.Range("M3").Value = 3.9149568210124843
$ws.Range("N3").Value = 3.450985516073886
$ws.Range("O3").Value = 3.4846402687314999
$ws.Range("P3").Value = 2.601772396032449
$ws.Range("Q3").Value = 2.2628828960183851
$ws.Range("R3").Value = 3.574492909694881
$ws.Range("S3").Value = 2.0699969550059567
$ws.Range("T3").Value = 3.0181062997856847
$ws.Range("U3").Value = 0.90158334739209689
$ws.Range("V3").Value = 0.90841098611064308
$ws.Range("W3").Value = 0.68496878593128507
$ws.Range("X3").Value = 0.20817752525871697

# Row 4
$ws.Range("C4").Value = -6.0382867568486684
$ws.Range("D4").Value = -6.5129714901912914
$ws.Range("E4").Value = 40.359746677123141
$ws.Range("F4").Value = 0.42506162981666107
$ws.Range("G4").Value = 1.4661518813452556
$ws.Range("H4").Value = 3.4682183225852263
$ws.Range("I4").Value = 12.03598930935105
$ws.Range("J4").Value = 4.1648031588382786
$ws.Range("K4").Value = 2.9484716317121284
$ws.Range("L4").Value = 3.4617709445197717
$ws.Range("M4").Value = 3.8307356378113382
$ws.Range("N4").Value = 1.4533504735832592
$ws.Range("O4").Value = 3.8113180359322367
$ws.Range("P4").Value = 2.9820199420242166
$ws.Range("Q4").Value = 2.7956340793634826
$ws.Range("R4").Value = 4.0905318878651498
$ws.Range("S4").Value = 2.366570101665662
$ws.Range("T4").Value = 3.4943762721235014
$ws.Range("U4").Value = 0.93748530750150627
$ws.Range("V4").Value = 0.95337751966252382
$ws.Range("W4").Value = 0.82497581497185257
$ws.Range("X4").Value = 0.28214594634898915

# Row 5
$ws.Range("C5").Value = -8.025338879250711
$ws.Range("D5").Value = -8.4549570367357632
$ws.Range("E5").Value = 59.087109840062261
$ws.Range("F5").Value = 0.77350611996255514
$ws.Range("G5").Value = 2.716537229351538
$ws.Range("H5").Value = 5.1631338718877746
$ws.Range("I5").Value = 8.2342659036933163
$ws.Range("J5").Value = 3.6560280516966834
$ws.Range("K5").Value = 2.2935280482498976
$ws.Range("L5").Value = 2.8615516141641204
$ws.Range("M5").Value = 3.4292117573950893
$ws.Range("N5").Value = 2.8938529170377203
$ws.Range("O5").Value = 2.8330499803292888
$ws.Range("P5").Value = 2.4443287998607275
$ws.Range("Q5").Value = 2.0687595282970013
$ws.Range("R5").Value = 3.223107946956616
$ws.Range("S5").Value = 1.8561171041386388
$ws.Range("T5").Value = 2.7440286318500084
$ws.Range("U5").Value = 0.77316881333802945
$ws.Range("V5").Value = 0.81599283138694234
$ws.Range("W5").Value = 0.57099320667547016
$ws.Range("X5").Value = 0.16026574709830702

